$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.901.72'
$ws.Range('E2').Value = '  -2.30%  '
$ws.Range('D3').Value = '1.798.22'
$ws.Range('E3').Value = '  -2.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.66'
$ws.Range('E5').Value = '  -7.30%  '
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5062'
$ws.Range('E7').Value = '  -3.72%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2510'
$ws.Range('E8').Value = '  -21.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06105'
$ws.Range('E9').Value = '  -10.17%  '
$ws.Range('D10').Value = '1.808.09'
$ws.Range('E10').Value = '  -1.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06874'
$ws.Range('E11').Value = '  -11.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.78'
$ws.Range('E12').Value = '  -21.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6140'
$ws.Range('E13').Value = '  -21.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '79.37'
$ws.Range('E14').Value = '  -9.74%  '
$ws.Range('E15').Value = '  -11.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.006'
$ws.Range('E17').Value = '  +0.52%  '
$ws.Range('D18').Value = '25.960.82'
$ws.Range('E18').Value = '  -2.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.36'
$ws.Range('E19').Value = '  -18.09%  '
$ws.Range('D20').Value = '2.052.61'
$ws.Range('E20').Value = '  -1.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006094'
$ws.Range('E21').Value = '  -23.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.959'
$ws.Range('E22').Value = '  -14.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.229'
$ws.Range('E23').Value = '  -12.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.068'
$ws.Range('E24').Value = '  -13.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '130.51'
$ws.Range('E25').Value = '  -8.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.887'
$ws.Range('E26').Value = '  -13.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.57'
$ws.Range('E27').Value = '  -13.82%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.299'
$ws.Range('E28').Value = '  -22.56%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '98.57'
$ws.Range('E29').Value = '  -11.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08259'
$ws.Range('E30').Value = '  -5.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.632'
$ws.Range('E31').Value = '  -12.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.752'
$ws.Range('E32').Value = '  -3.92%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04334'
$ws.Range('E33').Value = '  -11.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.164'
$ws.Range('E34').Value = '  -22.41%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.035'
$ws.Range('E35').Value = '  -8.88%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6320'
$ws.Range('E36').Value = '  -13.44%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.888'
$ws.Range('E37').Value = '  -6.70%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.099'
$ws.Range('E38').Value = '  -6.46%  '
$ws.Range('B39').Value = 'PaxDollar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.006'
$ws.Range('E39').Value = '  +0.52%  '
$ws.Range('B40').Value = 'Quant'
$ws.Range('C40').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '100.77'
$ws.Range('E40').Value = '  -8.10%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01469'
$ws.Range('E41').Value = '  -16.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7970'
$ws.Range('E42').Value = '  -10.93%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3911'
$ws.Range('E43').Value = '  -18.61%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.153'
$ws.Range('E44').Value = '  -13.35%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.217'
$ws.Range('E45').Value = '  -18.96%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05281'
$ws.Range('E46').Value = '  -9.72%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '53.46'
$ws.Range('E47').Value = '  -10.33%  '
$ws.Range('B48').Value = 'USDD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.007'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1034'
$ws.Range('E49').Value = '  -16.27%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '29.32'
$ws.Range('E50').Value = '  -16.04%  '
$ws.Range('B51').Value = 'TrueUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.002'
$ws.Range('E51').Value = '  +0.24%  '
